# Auto-generated edit script: refreshes crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.412.42'
$ws.Range("E2").Value = '  -2.04%  '

$ws.Range("D3").Value = '3.493.00'
$ws.Range("E3").Value = '  -2.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.09'
$ws.Range("E5").Value = '  +5.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.26'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.214'
$ws.Range("E9").Value = '  -2.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.652'
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.32'
$ws.Range("E11").Value = '  -2.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000308'
$ws.Range("E12").Value = '  -3.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.58'
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").Value = '4.049.65'
$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '605.98'
$ws.Range("E15").Value = '  +4.02%  '

$ws.Range("D16").Value = '69.402.39'
$ws.Range("E16").Value = '  -2.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.64'
$ws.Range("E17").Value = '  +1.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.85'
$ws.Range("E18").Value = '  -2.28%  '

$ws.Range("D19").Value = '3.501.25'
$ws.Range("E19").Value = '  -2.59%  '

$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.989'
$ws.Range("E21").Value = '  -1.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.28'
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.17'
$ws.Range("E23").Value = '  +9.69%  '

$ws.Range("E24").Value = '  +2.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.04'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.02'
$ws.Range("E26").Value = '  +2.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.96'
$ws.Range("E27").Value = '  -2.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  +8.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.68'
$ws.Range("E29").Value = '  +3.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.00'
$ws.Range("E30").Value = '  -3.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.46'
$ws.Range("E31").Value = '  +1.26%  '

$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("E33").Value = '  +16.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.32'
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("E35").Value = '  -6.70%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '522.77'
$ws.Range("E37").Value = '  -4.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.397'
$ws.Range("E38").Value = '  -4.55%  '

$ws.Range("D39").Value = '3.577.57'
$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.59'
$ws.Range("E40").Value = '  +5.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.77'
$ws.Range("E41").Value = '  -2.87%  '

$ws.Range("D42").Value = '0.0₃0775'
$ws.Range("E42").Value = '  -3.35%  '

$ws.Range("E43").Value = '  +0.96%  '

$ws.Range("E44").Value = '  +3.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.97'
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("E46").Value = '  +4.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("E47").Value = '  -5.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.84'
$ws.Range("E48").Value = '  -5.45%  '

$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.21'
$ws.Range("E50").Value = '  -3.48%  '

$ws.Range("E51").Value = '  -8.17%  '

